$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# Update N7 (Carry Look-Ahead Adder, Total Power column)
$ws.Range("N7").Value = 28.81

# Fill in Floating Point Adder dependency data on row 17
$ws.Range("Q17").Value = 1144.0
$ws.Range("R17").Value = 3758.4
$ws.Range("S17").Value = 16700.4
$ws.Range("T17").Value = 14741.6
$ws.Range("U17").Value = 99.17
$ws.Range("V17").Value = 3758.4
